$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The table shrinks from A1:I8 to A1:I5 -> drop the last three data rows.
$ws.Range("A6:I8").Delete()

# Columns B:H are always text-typed in this table (even numeric-looking values
# like "1", "360" stay as text), so force text format before writing so the
# COM layer doesn't silently coerce them to numbers.
$ws.Range("B2:H5").NumberFormat = "@"

# Row 2 (id 837 -> 393)
$ws.Range("A2").Value = 393
$ws.Range("B2").Value = "MC"
$ws.Range("C2").Value = "1"
$ws.Range("D2").Value = "no"
$ws.Range("E2").Value = "n.a."
$ws.Range("F2").Value = "n.a."
$ws.Range("G2").Value = "n.a."
$ws.Range("H2").Value = "yes"
$ws.Range("I2").Value = 11

# Row 3 (id 840 -> 398)
$ws.Range("A3").Value = 398
$ws.Range("B3").Value = "1"
$ws.Range("C3").Value = "1p"
$ws.Range("D3").Value = "yes"
$ws.Range("E3").Value = "360"
$ws.Range("F3").Value = "630"
$ws.Range("G3").Value = "30"
$ws.Range("H3").Value = "no"
$ws.Range("I3").Value = "working"

# Row 4 (id 845 -> 428)
$ws.Range("A4").Value = 428
$ws.Range("B4").Value = "1p"
$ws.Range("C4").Value = "MCd"
$ws.Range("D4").Value = "no"
$ws.Range("E4").Value = "n.a."
$ws.Range("F4").Value = "n.a."
$ws.Range("G4").Value = "n.a."
$ws.Range("H4").Value = "yes"
$ws.Range("I4").Value = 14

# Row 5 (id 875 -> 700)
$ws.Range("A5").Value = 700
$ws.Range("B5").Value = "MCd"
$ws.Range("C5").Value = "n.a."
$ws.Range("D5").Value = "yes"
$ws.Range("E5").Value = "n.a."
$ws.Range("F5").Value = "n.a."
$ws.Range("G5").Value = "n.a."
$ws.Range("H5").Value = "no"
$ws.Range("I5").Value = 14

# Writing through NumberFormat="@" bumps the cell style to a new text-format
# style; paste the original style (s=3, centered, thin border) back over the
# block from an untouched cell in the same style family so formatting is
# unchanged, as in the source diff.
$ws.Range("I2").Copy()
$ws.Range("B2:H5").PasteSpecial(-4122)
$ws.Range("I2").Value = 11
$excel.CutCopyMode = $false
